$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3: kamil's employee info (same employer "Walee"/"abc" as row 2)
$ws.Range("A3").Value = "kamil"
$ws.Range("B3").Value = "abc"
$ws.Range("C3").Value = "Walee"
$ws.Range("D3").Value = "abc"
$ws.Range("E3").Value = 111222333

# New row 4: Shah's employee info
$ws.Range("A4").Value = "Shah"
$ws.Range("B4").Value = "abc"
$ws.Range("C4").Value = "Walee"
$ws.Range("D4").Value = "abc "
$ws.Range("E4").Value = 222555888

# New columns F (doc_options) and G (paystub_options) with header + row2 values
$ws.Range("F1").Value = "doc_options"
$ws.Range("F2").Value = "1,2,3"
$ws.Range("G1").Value = "paystub_options"
$ws.Range("G2").Value = '{"Rate" : "", "Numbe of Paystubs":"", "Period" : "",  }'

# Widen the new doc_options column
$ws.Columns.Item(6).ColumnWidth = 11

# Ensure portrait page orientation (matches exported pageSetup)
$ws.PageSetup.Orientation = 1

# Leave selection on the last touched cell
$ws.Range("G4").Select()
